$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray G/H column cells (row 42) that are outside the real data range
$ws.Range("G1:H59").Clear()

# Apply integer number format + top vertical alignment across the whole data range
$rng = $ws.Range("A1:F59")
$rng.NumberFormat = "0"
$rng.VerticalAlignment = -4160

# Row 5 got manually resized
$ws.Rows("5").RowHeight = 19.5

# Update the selection to match what was left active in the saved file
$ws.Range("H5").Select()

Write-Output "done"
